# Applies the commit's cell-level edits across all four worksheets:
#   展览    (sheet 1) - F-column ("想去人数") value refreshes
#   演出    (sheet 2) - F-column ("想去人数") value refreshes
#   本地生活 (sheet 3) - row 4 (2024-04-30 entry) removed, rows shift up,
#                        A-column index renumbered, F-column values refreshed
#   全部类型 (sheet 4) - F-column ("想去人数") value refreshes

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1Changes = @(
    @(2, 1837),
    @(3, 20),
    @(5, 42),
    @(6, 1087),
    @(8, 180),
    @(9, 593),
    @(11, 456),
    @(12, 538),
    @(13, 1406),
    @(15, 1425),
    @(16, 21),
    @(17, 1173),
    @(20, 440),
    @(22, 315),
    @(25, 1275),
    @(26, 312),
    @(30, 1028),
    @(31, 0),
    @(32, 952),
    @(34, 1314),
    @(35, 882),
    @(36, 1040),
    @(37, 25),
    @(39, 1563),
    @(41, 27),
    @(42, 780),
    @(44, 767)
)
foreach ($pair in $ws1Changes) {
    $ws1.Cells.Item($pair[0], 6).Value = $pair[1]
}

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2Changes = @(
    @(4, 121),
    @(6, 140),
    @(10, 159),
    @(11, 1405),
    @(14, 2494),
    @(15, 1174),
    @(16, 387),
    @(18, 209),
    @(23, 423),
    @(27, 0),
    @(28, 14),
    @(31, 173),
    @(34, 55),
    @(38, 4),
    @(49, 6)
)
foreach ($pair in $ws2Changes) {
    $ws2.Cells.Item($pair[0], 6).Value = $pair[1]
}

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# The "2024-04-30" entry (old row 4) was dropped entirely; every row
# below it shifts up by one. Excel's own Rows.Delete() performs that
# shift for us (and keeps the A1:I16 -> A1:I15 used-range in sync), but
# it leaves the literal index numbers in column A untouched, so those
# are corrected by hand afterwards, along with the independent
# "想去人数" (F column) refreshes from the new scrape.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows(4).Delete()

# Renumber column A (index column) for rows 4..15: A(n) = n - 1
for ($r = 4; $r -le 15; $r++) {
    $ws3.Cells.Item($r, 1).Value = $r - 1
}

$ws3Changes = @(
    @(5, 2683),
    @(6, 4490),
    @(7, 116),
    @(9, 522),
    @(10, 619),
    @(11, 408),
    @(12, 178),
    @(13, 663),
    @(15, 349)
)
foreach ($pair in $ws3Changes) {
    $ws3.Cells.Item($pair[0], 6).Value = $pair[1]
}

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4Changes = @(
    @(2, 1837),
    @(4, 20),
    @(5, 4490),
    @(6, 619),
    @(7, 42),
    @(8, 178),
    @(9, 178),
    @(10, 663),
    @(11, 663),
    @(13, 140),
    @(14, 1087),
    @(16, 180),
    @(17, 159),
    @(18, 1405),
    @(19, 593),
    @(20, 456),
    @(21, 538),
    @(22, 2494),
    @(23, 1174),
    @(24, 1406),
    @(26, 1425),
    @(27, 1173),
    @(28, 209),
    @(32, 315),
    @(33, 349),
    @(34, 423),
    @(35, 1275),
    @(38, 1028),
    @(40, 952),
    @(41, 882),
    @(42, 25),
    @(45, 1563),
    @(49, 780),
    @(51, 767)
)
foreach ($pair in $ws4Changes) {
    $ws4.Cells.Item($pair[0], 6).Value = $pair[1]
}

Write-Host "Applied edits to all 4 sheets."
